# Auto-generated edit script: updates crypto price/volume table (rows 2-51)
# D-column (Price) values are forced to Text to preserve exact formatting
# (matches original string-typed cells; avoids Excel auto-numeric coercion
# that would drop trailing zeros / separators like "0.610" -> 0.61).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.462.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.286.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.65%  "

$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.618.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.870"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.298.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.388.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("E25").Value = "  +4.41%  "

$ws.Range("E26").Value = "  +1.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.29%  "

$ws.Range("E30").Value = "  -2.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("E34").Value = "  -1.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.38%  "

$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0368"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.62%  "

$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.78%  "

$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.94%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("E50").Value = "  +2.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
